# Auto-generated: apply Hortaliza/Vega Modelo de Temuco - Cebollín weekly update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 190
$ws.Cells.Item(190, 1).Value = 10
$ws.Cells.Item(190, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(190, 3).Value = 'La Araucanía'
$ws.Cells.Item(190, 4).Value = 44461
$ws.Cells.Item(190, 5).Value = 9
$ws.Cells.Item(190, 6).Value = 100112037
$ws.Cells.Item(190, 7).Value = 'Cebollín'
$ws.Cells.Item(190, 8).Value = 'Sin especificar'
$ws.Cells.Item(190, 9).Value = 'Primera'
$ws.Cells.Item(190, 10).Value = 20
$ws.Cells.Item(190, 11).Value = 8000
$ws.Cells.Item(190, 12).Value = 8000
$ws.Cells.Item(190, 13).Value = 8000
$ws.Cells.Item(190, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(190, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(190, 16).Value = 667
$ws.Cells.Item(190, 17).Value = 12
$ws.Cells.Item(190, 18).Value = 'Hortaliza'

# Row 191
$ws.Cells.Item(191, 1).Value = 10
$ws.Cells.Item(191, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(191, 3).Value = 'La Araucanía'
$ws.Cells.Item(191, 4).Value = 44461
$ws.Cells.Item(191, 5).Value = 9
$ws.Cells.Item(191, 6).Value = 100112037
$ws.Cells.Item(191, 7).Value = 'Cebollín'
$ws.Cells.Item(191, 8).Value = 'Sin especificar'
$ws.Cells.Item(191, 9).Value = 'Primera'
$ws.Cells.Item(191, 10).Value = 40
$ws.Cells.Item(191, 11).Value = 5000
$ws.Cells.Item(191, 12).Value = 5000
$ws.Cells.Item(191, 13).Value = 5000
$ws.Cells.Item(191, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(191, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(191, 16).Value = 417
$ws.Cells.Item(191, 17).Value = 12
$ws.Cells.Item(191, 18).Value = 'Hortaliza'

# Row 192
$ws.Cells.Item(192, 1).Value = 10
$ws.Cells.Item(192, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(192, 3).Value = 'La Araucanía'
$ws.Cells.Item(192, 4).Value = 44357
$ws.Cells.Item(192, 5).Value = 9
$ws.Cells.Item(192, 6).Value = 100112037
$ws.Cells.Item(192, 7).Value = 'Cebollín'
$ws.Cells.Item(192, 8).Value = 'Sin especificar'
$ws.Cells.Item(192, 9).Value = 'Primera'
$ws.Cells.Item(192, 10).Value = 80
$ws.Cells.Item(192, 11).Value = 9000
$ws.Cells.Item(192, 12).Value = 9000
$ws.Cells.Item(192, 13).Value = 9000
$ws.Cells.Item(192, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(192, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(192, 16).Value = 750
$ws.Cells.Item(192, 17).Value = 12
$ws.Cells.Item(192, 18).Value = 'Hortaliza'

# Row 193
$ws.Cells.Item(193, 1).Value = 10
$ws.Cells.Item(193, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(193, 3).Value = 'La Araucanía'
$ws.Cells.Item(193, 4).Value = 44203
$ws.Cells.Item(193, 5).Value = 9
$ws.Cells.Item(193, 6).Value = 100112037
$ws.Cells.Item(193, 7).Value = 'Cebollín'
$ws.Cells.Item(193, 8).Value = 'Sin especificar'
$ws.Cells.Item(193, 9).Value = 'Primera'
$ws.Cells.Item(193, 10).Value = 100
$ws.Cells.Item(193, 11).Value = 8000
$ws.Cells.Item(193, 12).Value = 8000
$ws.Cells.Item(193, 13).Value = 8000
$ws.Cells.Item(193, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(193, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(193, 16).Value = 667
$ws.Cells.Item(193, 17).Value = 12
$ws.Cells.Item(193, 18).Value = 'Hortaliza'

# Row 194
$ws.Cells.Item(194, 1).Value = 10
$ws.Cells.Item(194, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(194, 3).Value = 'La Araucanía'
$ws.Cells.Item(194, 4).Value = 44162
$ws.Cells.Item(194, 5).Value = 9
$ws.Cells.Item(194, 6).Value = 100112037
$ws.Cells.Item(194, 7).Value = 'Cebollín'
$ws.Cells.Item(194, 8).Value = 'Sin especificar'
$ws.Cells.Item(194, 9).Value = 'Primera'
$ws.Cells.Item(194, 10).Value = 50
$ws.Cells.Item(194, 11).Value = 7000
$ws.Cells.Item(194, 12).Value = 7000
$ws.Cells.Item(194, 13).Value = 7000
$ws.Cells.Item(194, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(194, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(194, 16).Value = 583
$ws.Cells.Item(194, 17).Value = 12
$ws.Cells.Item(194, 18).Value = 'Hortaliza'

# Row 195
$ws.Cells.Item(195, 1).Value = 10
$ws.Cells.Item(195, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(195, 3).Value = 'La Araucanía'
$ws.Cells.Item(195, 4).Value = 44162
$ws.Cells.Item(195, 5).Value = 9
$ws.Cells.Item(195, 6).Value = 100112037
$ws.Cells.Item(195, 7).Value = 'Cebollín'
$ws.Cells.Item(195, 8).Value = 'Sin especificar'
$ws.Cells.Item(195, 9).Value = 'Primera'
$ws.Cells.Item(195, 10).Value = 40
$ws.Cells.Item(195, 11).Value = 5000
$ws.Cells.Item(195, 12).Value = 5000
$ws.Cells.Item(195, 13).Value = 5000
$ws.Cells.Item(195, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(195, 15).Value = 'Región del Maule'
$ws.Cells.Item(195, 16).Value = 417
$ws.Cells.Item(195, 17).Value = 12
$ws.Cells.Item(195, 18).Value = 'Hortaliza'

# Row 196
$ws.Cells.Item(196, 1).Value = 10
$ws.Cells.Item(196, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(196, 3).Value = 'La Araucanía'
$ws.Cells.Item(196, 4).Value = 44410
$ws.Cells.Item(196, 5).Value = 9
$ws.Cells.Item(196, 6).Value = 100112037
$ws.Cells.Item(196, 7).Value = 'Cebollín'
$ws.Cells.Item(196, 8).Value = 'Sin especificar'
$ws.Cells.Item(196, 9).Value = 'Primera'
$ws.Cells.Item(196, 10).Value = 60
$ws.Cells.Item(196, 11).Value = 9000
$ws.Cells.Item(196, 12).Value = 10000
$ws.Cells.Item(196, 13).Value = 9500
$ws.Cells.Item(196, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(196, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(196, 16).Value = 792
$ws.Cells.Item(196, 17).Value = 12
$ws.Cells.Item(196, 18).Value = 'Hortaliza'

# Row 197
$ws.Cells.Item(197, 1).Value = 10
$ws.Cells.Item(197, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(197, 3).Value = 'La Araucanía'
$ws.Cells.Item(197, 4).Value = 44410
$ws.Cells.Item(197, 5).Value = 9
$ws.Cells.Item(197, 6).Value = 100112037
$ws.Cells.Item(197, 7).Value = 'Cebollín'
$ws.Cells.Item(197, 8).Value = 'Sin especificar'
$ws.Cells.Item(197, 9).Value = 'Primera'
$ws.Cells.Item(197, 10).Value = 50
$ws.Cells.Item(197, 11).Value = 5000
$ws.Cells.Item(197, 12).Value = 5000
$ws.Cells.Item(197, 13).Value = 5000
$ws.Cells.Item(197, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(197, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(197, 16).Value = 417
$ws.Cells.Item(197, 17).Value = 12
$ws.Cells.Item(197, 18).Value = 'Hortaliza'

# Row 198
$ws.Cells.Item(198, 1).Value = 10
$ws.Cells.Item(198, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(198, 3).Value = 'La Araucanía'
$ws.Cells.Item(198, 4).Value = 44411
$ws.Cells.Item(198, 5).Value = 9
$ws.Cells.Item(198, 6).Value = 100112037
$ws.Cells.Item(198, 7).Value = 'Cebollín'
$ws.Cells.Item(198, 8).Value = 'Sin especificar'
$ws.Cells.Item(198, 9).Value = 'Primera'
$ws.Cells.Item(198, 10).Value = 30
$ws.Cells.Item(198, 11).Value = 10000
$ws.Cells.Item(198, 12).Value = 10000
$ws.Cells.Item(198, 13).Value = 10000
$ws.Cells.Item(198, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(198, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(198, 16).Value = 833
$ws.Cells.Item(198, 17).Value = 12
$ws.Cells.Item(198, 18).Value = 'Hortaliza'

# Row 199
$ws.Cells.Item(199, 1).Value = 10
$ws.Cells.Item(199, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(199, 3).Value = 'La Araucanía'
$ws.Cells.Item(199, 4).Value = 44176
$ws.Cells.Item(199, 5).Value = 9
$ws.Cells.Item(199, 6).Value = 100112037
$ws.Cells.Item(199, 7).Value = 'Cebollín'
$ws.Cells.Item(199, 8).Value = 'Sin especificar'
$ws.Cells.Item(199, 9).Value = 'Primera'
$ws.Cells.Item(199, 10).Value = 30
$ws.Cells.Item(199, 11).Value = 8000
$ws.Cells.Item(199, 12).Value = 9000
$ws.Cells.Item(199, 13).Value = 8667
$ws.Cells.Item(199, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(199, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(199, 16).Value = 722
$ws.Cells.Item(199, 17).Value = 12
$ws.Cells.Item(199, 18).Value = 'Hortaliza'

# Row 200
$ws.Cells.Item(200, 1).Value = 10
$ws.Cells.Item(200, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(200, 3).Value = 'La Araucanía'
$ws.Cells.Item(200, 4).Value = 44239
$ws.Cells.Item(200, 5).Value = 9
$ws.Cells.Item(200, 6).Value = 100112037
$ws.Cells.Item(200, 7).Value = 'Cebollín'
$ws.Cells.Item(200, 8).Value = 'Sin especificar'
$ws.Cells.Item(200, 9).Value = 'Primera'
$ws.Cells.Item(200, 10).Value = 115
$ws.Cells.Item(200, 11).Value = 6000
$ws.Cells.Item(200, 12).Value = 7000
$ws.Cells.Item(200, 13).Value = 6565
$ws.Cells.Item(200, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(200, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(200, 16).Value = 547
$ws.Cells.Item(200, 17).Value = 12
$ws.Cells.Item(200, 18).Value = 'Hortaliza'

# Row 201
$ws.Cells.Item(201, 1).Value = 10
$ws.Cells.Item(201, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(201, 3).Value = 'La Araucanía'
$ws.Cells.Item(201, 4).Value = 44376
$ws.Cells.Item(201, 5).Value = 9
$ws.Cells.Item(201, 6).Value = 100112037
$ws.Cells.Item(201, 7).Value = 'Cebollín'
$ws.Cells.Item(201, 8).Value = 'Sin especificar'
$ws.Cells.Item(201, 9).Value = 'Primera'
$ws.Cells.Item(201, 10).Value = 65
$ws.Cells.Item(201, 11).Value = 5000
$ws.Cells.Item(201, 12).Value = 5000
$ws.Cells.Item(201, 13).Value = 5000
$ws.Cells.Item(201, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(201, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(201, 16).Value = 417
$ws.Cells.Item(201, 17).Value = 12
$ws.Cells.Item(201, 18).Value = 'Hortaliza'

# Row 202
$ws.Cells.Item(202, 1).Value = 10
$ws.Cells.Item(202, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(202, 3).Value = 'La Araucanía'
$ws.Cells.Item(202, 4).Value = 44292
$ws.Cells.Item(202, 5).Value = 9
$ws.Cells.Item(202, 6).Value = 100112037
$ws.Cells.Item(202, 7).Value = 'Cebollín'
$ws.Cells.Item(202, 8).Value = 'Sin especificar'
$ws.Cells.Item(202, 9).Value = 'Primera'
$ws.Cells.Item(202, 10).Value = 35
$ws.Cells.Item(202, 11).Value = 7000
$ws.Cells.Item(202, 12).Value = 7000
$ws.Cells.Item(202, 13).Value = 7000
$ws.Cells.Item(202, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(202, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(202, 16).Value = 583
$ws.Cells.Item(202, 17).Value = 12
$ws.Cells.Item(202, 18).Value = 'Hortaliza'

# Row 203
$ws.Cells.Item(203, 1).Value = 10
$ws.Cells.Item(203, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(203, 3).Value = 'La Araucanía'
$ws.Cells.Item(203, 4).Value = 44358
$ws.Cells.Item(203, 5).Value = 9
$ws.Cells.Item(203, 6).Value = 100112037
$ws.Cells.Item(203, 7).Value = 'Cebollín'
$ws.Cells.Item(203, 8).Value = 'Sin especificar'
$ws.Cells.Item(203, 9).Value = 'Primera'
$ws.Cells.Item(203, 10).Value = 30
$ws.Cells.Item(203, 11).Value = 9000
$ws.Cells.Item(203, 12).Value = 9000
$ws.Cells.Item(203, 13).Value = 9000
$ws.Cells.Item(203, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(203, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(203, 16).Value = 750
$ws.Cells.Item(203, 17).Value = 12
$ws.Cells.Item(203, 18).Value = 'Hortaliza'

# Row 204
$ws.Cells.Item(204, 1).Value = 10
$ws.Cells.Item(204, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(204, 3).Value = 'La Araucanía'
$ws.Cells.Item(204, 4).Value = 44211
$ws.Cells.Item(204, 5).Value = 9
$ws.Cells.Item(204, 6).Value = 100112037
$ws.Cells.Item(204, 7).Value = 'Cebollín'
$ws.Cells.Item(204, 8).Value = 'Sin especificar'
$ws.Cells.Item(204, 9).Value = 'Primera'
$ws.Cells.Item(204, 10).Value = 110
$ws.Cells.Item(204, 11).Value = 8000
$ws.Cells.Item(204, 12).Value = 8000
$ws.Cells.Item(204, 13).Value = 8000
$ws.Cells.Item(204, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(204, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(204, 16).Value = 667
$ws.Cells.Item(204, 17).Value = 12
$ws.Cells.Item(204, 18).Value = 'Hortaliza'

# Row 205
$ws.Cells.Item(205, 1).Value = 10
$ws.Cells.Item(205, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(205, 3).Value = 'La Araucanía'
$ws.Cells.Item(205, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(205, 4).Value = 44425
$ws.Cells.Item(205, 5).Value = 9
$ws.Cells.Item(205, 6).Value = 100112037
$ws.Cells.Item(205, 7).Value = 'Cebollín'
$ws.Cells.Item(205, 8).Value = 'Sin especificar'
$ws.Cells.Item(205, 9).Value = 'Primera'
$ws.Cells.Item(205, 10).Value = 30
$ws.Cells.Item(205, 11).Value = 7000
$ws.Cells.Item(205, 12).Value = 7000
$ws.Cells.Item(205, 13).Value = 7000
$ws.Cells.Item(205, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(205, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(205, 16).Value = 583
$ws.Cells.Item(205, 17).Value = 12
$ws.Cells.Item(205, 18).Value = 'Hortaliza'

# Row 206
$ws.Cells.Item(206, 1).Value = 10
$ws.Cells.Item(206, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(206, 3).Value = 'La Araucanía'
$ws.Cells.Item(206, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(206, 4).Value = 44425
$ws.Cells.Item(206, 5).Value = 9
$ws.Cells.Item(206, 6).Value = 100112037
$ws.Cells.Item(206, 7).Value = 'Cebollín'
$ws.Cells.Item(206, 8).Value = 'Sin especificar'
$ws.Cells.Item(206, 9).Value = 'Primera'
$ws.Cells.Item(206, 10).Value = 20
$ws.Cells.Item(206, 11).Value = 5000
$ws.Cells.Item(206, 12).Value = 5000
$ws.Cells.Item(206, 13).Value = 5000
$ws.Cells.Item(206, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(206, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(206, 16).Value = 417
$ws.Cells.Item(206, 17).Value = 12
$ws.Cells.Item(206, 18).Value = 'Hortaliza'
